# Apply "Attempt to add 7th" edit to the Chords sheet.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Basic Notes")
$ws2 = $wb.Worksheets.Item("Chords")

# --- Add a new "N" marker into C, D, E for every data row (2-43) ---
for ($r = 2; $r -le 43; $r++) {
    $ws2.Range("C$r").Value = "N"
    $ws2.Range("D$r").Value = "N"
    $ws2.Range("E$r").Value = "N"
}

# --- Turn header cell C1 into text "7th" (quote-prefixed, like typing '7th) ---
$ws2.Range("C1").Formula = "'7th"

# --- Apply an AutoFilter over the whole table ---
$tableRng = $ws2.Range("A1:G43")
$tableRng.AutoFilter()

# --- Sort the table by column A (Note), ascending, same as the author's attempt ---
$ws2.Sort.SortFields.Clear()
$ws2.Sort.SortFields.Add($ws2.Range("A1:A43"))
$ws2.Sort.SetRange($tableRng)
$ws2.Sort.Header = 1
$ws2.Sort.Apply()

# --- Register the hidden _FilterDatabase name Excel creates for AutoFilter+Sort ---
$filterName = $ws2.Names.Add('_xlnm._FilterDatabase', '=Chords!$A$1:$G$43')
$filterName.Visible = $false

# --- Make "Chords" the active sheet / tab, with the selection left on rows 8:11 ---
$ws2.Activate()
$ws2.Range("A8:XFD11").Select()
